$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that Excel would otherwise "smart-convert" (e.g. a
# date-like string such as "11.4.21") as plain text, preserving the
# destination cell's existing style. We build the text in a scratch cell via
# a formula (so it evaluates to a pure string result), copy it, and paste
# "values only" into the destination - this brings over the text without
# touching the destination's number format/style.
function Set-TextValue {
    param($cell, [string]$text)

    $scratch = $ws.Range("ZZ1")
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $ws.Range($cell).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
    $scratch.Clear()
    $excel.CutCopyMode = $false
}

# --- New rows for subjects 1010 and 1011 (11.4.21 full-arm-reach timing runs) ---

Set-TextValue "A23" "11.4.21"
$ws.Range("B23").Value = 1010
$ws.Range("P23").Value = "118B"
$ws.Range("T23").Value = "Khen heller"
Set-TextValue "U23" "Full run on myself with full arm reach to check to make sure all ok before running subs"
$ws.Rows(23).RowHeight = 15.75

Set-TextValue "A24" "11.4.21"
$ws.Range("B24").Value = 1011
$ws.Range("P24").Value = "118B"
$ws.Range("T24").Value = "Khen heller"
Set-TextValue "U24" "Full run on myself with full arm reach to check timing on diff screen (Asus, refrate 100) "
$ws.Rows(24).RowHeight = 15.75

# --- Fill in remaining demographic fields for subject #2 (row 4) ---

$ws.Range("I4").Value = "M"
Set-TextValue "A4" "12.4.21"
$ws.Range("C4").Value = 30
$ws.Range("F4").Value = "right"
$ws.Range("G4").Value = "N"
$ws.Range("J4").Value = "Y"
$ws.Range("L4").Value = "N"
$ws.Range("M4").Value = "N"
$ws.Range("N4").Value = "N"

# --- Restore the view: scroll back to column A and select J9 ---
$ws.Range("J9").Select()
